$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.619.52"
$ws.Range("E2").Value = "  +2.46%  "
$ws.Range("D3").Value = "2.410.39"
$ws.Range("E3").Value = "  +3.01%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'552.72"
$ws.Range("E5").Value = "  +2.16%  "
$ws.Range("D6").Value = "'136.90"
$ws.Range("E6").Value = "  +1.57%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E9").Value = "  +5.36%  "
$ws.Range("D10").Value = "'5.80"
$ws.Range("E10").Value = "  +2.80%  "
$ws.Range("E11").Value = "  +1.86%  "
$ws.Range("E12").Value = "  -2.08%  "
$ws.Range("D14").Value = "2.837.38"
$ws.Range("E14").Value = "  +2.93%  "
$ws.Range("D15").Value = "59.501.45"
$ws.Range("E15").Value = "  +2.35%  "
$ws.Range("D16").Value = "'0.0000140"
$ws.Range("E16").Value = "  +4.50%  "
$ws.Range("D17").Value = "2.411.36"
$ws.Range("E17").Value = "  +2.78%  "
$ws.Range("E18").Value = "  +5.84%  "
$ws.Range("D19").Value = "'4.46"
$ws.Range("E19").Value = "  +4.84%  "
$ws.Range("D20").Value = "'335.91"
$ws.Range("E20").Value = "  +0.87%  "
$ws.Range("E21").Value = "  +4.74%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").Value = "'64.71"
$ws.Range("E23").Value = "  +3.14%  "
$ws.Range("E24").Value = "  +0.80%  "
$ws.Range("D25").Value = "'8.49"
$ws.Range("E25").Value = "  -0.54%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("E27").Value = "  -1.96%  "
$ws.Range("D28").Value = "0.0₃0781"
$ws.Range("E28").Value = "  +6.30%  "
$ws.Range("E29").Value = "  +2.25%  "
$ws.Range("D30").Value = "'170.68"
$ws.Range("E30").Value = "  +0.27%  "
$ws.Range("D31").Value = "'6.25"
$ws.Range("E31").Value = "  +2.50%  "
$ws.Range("E32").Value = "  +1.62%  "
$ws.Range("D33").Value = "'1.03"
$ws.Range("E33").Value = "  -0.29%  "
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'1.31"
$ws.Range("E35").Value = "  +5.90%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").Value = "'4.31"
$ws.Range("E36").Value = "  +1.73%  "
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("E38").Value = "  +0.17%  "
$ws.Range("D39").Value = "'40.09"
$ws.Range("E39").Value = "  +2.65%  "
$ws.Range("D40").Value = "'0.418"
$ws.Range("E40").Value = "  +11.23%  "
$ws.Range("D41").Value = "'304.16"
$ws.Range("E41").Value = "  +6.31%  "
$ws.Range("E42").Value = "  +3.03%  "
$ws.Range("D43").Value = "'142.44"
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("E44").Value = "  +2.76%  "
$ws.Range("E45").Value = "  +4.71%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Value = "'0.571"
$ws.Range("E46").Value = "  +1.37%  "
$ws.Range("B47").Value = "Polygon"
$ws.Range("C47").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D47").Value = "'0.407"
$ws.Range("E47").Value = "  +6.59%  "
$ws.Range("D48").Value = "'19.06"
$ws.Range("E48").Value = "  -0.47%  "
$ws.Range("E49").Value = "  +3.53%  "
$ws.Range("D50").Value = "'11.04"
$ws.Range("E50").Value = "  -0.29%  "
$ws.Range("D51").Value = "'1.60"
$ws.Range("E51").Value = "  +5.03%  "
